$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.944.66"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "3.088.96"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.31"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.24"
$ws.Range("E6").Value = "  +6.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.084.73"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.33"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "3.603.68"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("D17").Value = "66.880.72"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "3.092.79"
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.22"
$ws.Range("E20").Value = "  +4.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.59"
$ws.Range("E21").Value = "  +4.61%  "
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.03"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.07"
$ws.Range("E25").Value = "  +7.38%  "
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.39"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +3.67%  "
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.67"
$ws.Range("E38").Value = "  +7.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.08"
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.55"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.65"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0359"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "380.61"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "2.783.80"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.91"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.75"
$ws.Range("E50").Value = "  +5.56%  "
$ws.Range("E51").Value = "  +0.90%  "
